$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C3" = -13.406
    "B7" = 5.697000000000001
    "A8" = -22.084
    "A10" = -21.598
    "A12" = -21.515
    "B15" = 5.072000000000001
    "A18" = -21.498
    "B18" = 6.989
    "C18" = -11.919
    "C19" = -11.706
    "B20" = 6.956999999999999
    "C27" = -13.103
    "B29" = 5.645
    "B30" = 6.208
    "B31" = 5.396
    "C31" = -13.066
    "A37" = -20.287
    "C38" = -12.97
    "B40" = 8.962
    "C42" = -12.533
    "C44" = -12.777
    "C47" = -12.661
    "B50" = 5.896000000000001
    "A55" = -21.826
    "C58" = -13.09
    "C65" = -12.282
    "A68" = -21.667
    "B68" = 5.422
    "C73" = -12.57
    "B76" = 6.689
    "A77" = -20.609
    "A78" = -20.083
    "A81" = -21.906
    "A82" = -21.795
    "B87" = 5.508999999999999
    "B88" = 5.816000000000001
    "C90" = -13.321
    "C94" = -10.505
    "C95" = -11.456
    "B96" = 7.187
    "B98" = 5.49
    "B101" = 7.840999999999999
    "C101" = -12.968
    "B102" = 7.306999999999999
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

